$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments
$ws.Columns.Item(1).ColumnWidth = 20.83
$ws.Columns.Item(19).ColumnWidth = 21.83

# Corrected values: uncorrected age calculation (col A) and isotope masses
# ripple through derived ratio/error columns (B, O, P, S, T)
$ws.Range("A2").Value = -1.371845040908592
$ws.Range("B2").Value = 0.001669604528143381
$ws.Range("O2").Value = 1.324582333196705
$ws.Range("P2").Value = 0.3409005173764483
$ws.Range("S2").Value = 0.2425649162228197
$ws.Range("T2").Value = 0.3269535019506545
$ws.Range("A3").Value = 144.3120226285348
$ws.Range("B3").Value = 0.0005475686589663281
$ws.Range("O3").Value = 0.003275356462361479
$ws.Range("P3").Value = 2.119924914268915
$ws.Range("S3").Value = 0.0003203543385384721
$ws.Range("T3").Value = 2.020137790100577
$ws.Range("A4").Value = -1.410909894590029
$ws.Range("B4").Value = 0.001844912777490458
$ws.Range("O4").Value = 1.32333141153102
$ws.Range("P4").Value = 0.2852828322387429
$ws.Range("S4").Value = 0.2414261194212968
$ws.Range("T4").Value = 0.3636864051450315
$ws.Range("A5").Value = 144.2537039446266
$ws.Range("B5").Value = 0.0004965239513688518
$ws.Range("O5").Value = 0.002811543931097559
$ws.Range("P5").Value = 1.998927788018716
$ws.Range("S5").Value = 0.0005106998173672897
$ws.Range("T5").Value = 2.009498076031131
$ws.Range("A6").Value = 1.024861744547856
$ws.Range("B6").Value = 0.001799480281915367
$ws.Range("O6").Value = 1.32219823094449
$ws.Range("P6").Value = 0.2209739728222805
$ws.Range("S6").Value = 0.2396666663667818
$ws.Range("T6").Value = 0.3656704898636055
$ws.Range("A7").Value = 144.6551671144802
$ws.Range("B7").Value = 0.0005465529888301912
$ws.Range("O7").Value = 0.003652908602329808
$ws.Range("P7").Value = 1.937390319586194
$ws.Range("S7").Value = 0.0006632840911626782
$ws.Range("T7").Value = 1.940759177824897
$ws.Range("A8").Value = -0.7891165434938952
$ws.Range("B8").Value = 0.001331415506142307
$ws.Range("O8").Value = 1.319622257679337
$ws.Range("P8").Value = 0.2325182917524557
$ws.Range("S8").Value = 0.2382819122323587
$ws.Range("T8").Value = 0.3670471312898169
$ws.Range("A9").Value = 144.2900492926518
$ws.Range("B9").Value = 0.0003788151492977407
$ws.Range("O9").Value = 0.003179041188988946
$ws.Range("P9").Value = 2.1377126504212
$ws.Range("S9").Value = 0.0003630052528144143
$ws.Range("T9").Value = 2.137402394568339
$ws.Range("A10").Value = -0.7824415973383125
$ws.Range("B10").Value = 0.001926583093367701
$ws.Range("O10").Value = 1.317299619129553
$ws.Range("P10").Value = 0.2453970110927369
$ws.Range("S10").Value = 0.2367852110607417
$ws.Range("T10").Value = 0.4260308980406651
$ws.Range("A11").Value = 145.2560364805913
$ws.Range("B11").Value = 0.0004874081370088768
$ws.Range("O11").Value = 0.003396968616918185
$ws.Range("P11").Value = 2.26423125468954
$ws.Range("S11").Value = 0.0007049510064866325
$ws.Range("T11").Value = 2.297564997300947
$ws.Range("A12").Value = -0.3666892951877943
$ws.Range("B12").Value = 0.001503954364435659
$ws.Range("O12").Value = 1.316499874313986
$ws.Range("P12").Value = 0.2613402682972835
$ws.Range("S12").Value = 0.2383381790433829
$ws.Range("T12").Value = 0.4768270850367027
